$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column B width (target OOXML width="11.7109375").
# The host's ColumnWidth setter quantizes to 1/6-character steps, so the
# closest attainable OOXML width is 11.666666666666666 (ColumnWidth=65/6).
$ws.Columns.Item(2).ColumnWidth = 10.833333333333334

# Update cell values
$ws.Range("A1").Value = 148.84975220278793
$ws.Range("B1").Value = 4.6454898084276053
$ws.Range("C1").Value = 0.7670254403131116
